# Refresh the cryptocurrency price/volume snapshot (columns D and E,
# rows 2-51) to the latest scraped values.
#
# Column D holds price text. Some of the new prices are plain numeric
# strings (e.g. "0.9994"), and Excel auto-converts those to real numbers
# on a normal .Value assignment, which would change the stored cell type
# away from text. For those rows we prefix the input with an apostrophe
# (exactly like a user typing '0.9994 into Excel) to force a text entry,
# then reset the cell Style back to "Normal" so the quote-prefix cell
# format it creates is not left behind. Prices that already parse as
# non-numeric text (e.g. "30.436.51", with two dots) need no such fixup.
#
# Column E holds the padded percentage-change strings; the surrounding
# spaces already keep Excel from treating them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "30.436.51"; DAsText = $false; E = "  -1.24%  " }
    @{ Row = 3; D = "1.915.37"; DAsText = $false; E = "  +1.50%  " }
    @{ Row = 4; D = "0.9994"; DAsText = $true; E = "  -0.10%  " }
    @{ Row = 5; D = "242.77"; DAsText = $true; E = "  +1.75%  " }
    @{ Row = 6; D = "0.9998"; DAsText = $true; E = "  -0.04%  " }
    @{ Row = 7; D = "0.4694"; DAsText = $true; E = "  -1.68%  " }
    @{ Row = 8; D = "0.2863"; DAsText = $true; E = "  -0.86%  " }
    @{ Row = 9; D = "0.06813"; DAsText = $true; E = "  +3.48%  " }
    @{ Row = 10; D = "109.97"; DAsText = $true; E = "  +13.01%  " }
    @{ Row = 11; D = "18.46"; DAsText = $true; E = "  -2.36%  " }
    @{ Row = 12; D = "0.07723"; DAsText = $true; E = "  +1.66%  " }
    @{ Row = 13; D = "1.896.47"; DAsText = $false; E = "  +0.57%  " }
    @{ Row = 14; D = "5.279"; DAsText = $true; E = "  +3.05%  " }
    @{ Row = 15; D = "0.6585"; DAsText = $true; E = "  +0.04%  " }
    @{ Row = 16; D = "295.84"; DAsText = $true; E = "  -3.54%  " }
    @{ Row = 17; D = "30.423.70"; DAsText = $false; E = "  -1.24%  " }
    @{ Row = 18; D = "0.000007619"; DAsText = $true; E = "  +0.44%  " }
    @{ Row = 19; D = "1.000"; DAsText = $true; E = "  -0.06%  " }
    @{ Row = 20; D = "12.91"; DAsText = $true; E = "  -2.11%  " }
    @{ Row = 21; D = "2.132.34"; DAsText = $false; E = "  +0.36%  " }
    @{ Row = 22; D = "0.9985"; DAsText = $true; E = "  -0.10%  " }
    @{ Row = 23; D = "5.244"; DAsText = $true; E = "  +2.32%  " }
    @{ Row = 24; D = "6.202"; DAsText = $true; E = "  +0.64%  " }
    @{ Row = 25; D = "21.78"; DAsText = $true; E = "  +6.95%  " }
    @{ Row = 26; D = "9.315"; DAsText = $true; E = "  +0.18%  " }
    @{ Row = 27; D = "168.33"; DAsText = $true; E = "  +1.15%  " }
    @{ Row = 28; D = "2.087"; DAsText = $true; E = "  +7.09%  " }
    @{ Row = 29; D = $null; DAsText = $false; E = "  +0.00%  " }
    @{ Row = 30; D = "1.365"; DAsText = $true; E = "  +0.78%  " }
    @{ Row = 31; D = "4.179"; DAsText = $true; E = "  +0.07%  " }
    @{ Row = 32; D = "3.989"; DAsText = $true; E = "  +0.39%  " }
    @{ Row = 33; D = "0.05054"; DAsText = $true; E = "  +0.33%  " }
    @{ Row = 34; D = "0.7376"; DAsText = $true; E = "  +0.80%  " }
    @{ Row = 35; D = "1.154"; DAsText = $true; E = "  -1.39%  " }
    @{ Row = 36; D = "0.02072"; DAsText = $true; E = "  +6.50%  " }
    @{ Row = 37; D = $null; DAsText = $false; E = "  +1.14%  " }
    @{ Row = 38; D = "2.675"; DAsText = $true; E = "  -0.90%  " }
    @{ Row = 39; D = $null; DAsText = $false; E = "  -1.31%  " }
    @{ Row = 40; D = "109.28"; DAsText = $true; E = "  +1.44%  " }
    @{ Row = 41; D = "0.8704"; DAsText = $true; E = "  -3.83%  " }
    @{ Row = 42; D = "5.820"; DAsText = $true; E = "  +3.07%  " }
    @{ Row = 43; D = "0.4265"; DAsText = $true; E = "  +1.45%  " }
    @{ Row = 44; D = "0.9995"; DAsText = $true; E = "  -0.09%  " }
    @{ Row = 45; D = "51.55"; DAsText = $true; E = "  +20.49%  " }
    @{ Row = 46; D = "67.48"; DAsText = $true; E = "  +2.69%  " }
    @{ Row = 47; D = $null; DAsText = $false; E = "  -2.25%  " }
    @{ Row = 48; D = "9.249"; DAsText = $true; E = "  +2.70%  " }
    @{ Row = 49; D = "0.1220"; DAsText = $true; E = "  -0.37%  " }
    @{ Row = 50; D = "34.88"; DAsText = $true; E = "  +0.22%  " }
    @{ Row = 51; D = "0.2456"; DAsText = $true; E = "  +11.69%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Range("D" + $u.Row)
        if ($u.DAsText) {
            $cell.Value = "'" + $u.D
            $cell.Style = "Normal"
        } else {
            $cell.Value = $u.D
        }
    }
    if ($null -ne $u.E) {
        $ws.Range("E" + $u.Row).Value = $u.E
    }
}
